$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The last existing data row is row 135 (date serial 45691).
# Append two new rows (136, 137) that replicate the same values as row 135,
# only advancing the date in column A by one day each time.
$lastRow = 135

for ($i = 1; $i -le 2; $i++) {
    $srcRow = $lastRow
    $dstRow = $lastRow + $i

    # Copy the whole row (values, number formats, styles) down to the new row.
    $ws.Range("A$srcRow`:J$srcRow").Copy($ws.Range("A$dstRow`:J$dstRow"))

    # Update the date serial in column A for the new row.
    $ws.Range("A$dstRow").Value = 45691 + $i
}
